# Update Polinago.xlsx with new daily rows through 23 August 2021
# (aggiornamento al 23 agosto 2021)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Last existing data row is 343 (date serial 44417 = 2021-08-09).
# Use it as a style/format template for the appended rows.
$lastRow = 343
$template = $ws.Range("A$($lastRow):D$($lastRow)")

# New rows to append: date serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila ab.
$newData = @(
    @(44418, 0, 0, 0),
    @(44419, 0, 0, 0),
    @(44420, 2, 2, 124.6882793017456),
    @(44421, 0, 2, 124.6882793017456),
    @(44422, 0, 2, 124.6882793017456),
    @(44423, 1, 3, 187.0324189526185),
    @(44424, 0, 3, 187.0324189526185),
    @(44425, 0, 3, 187.0324189526185),
    @(44426, 2, 5, 311.7206982543641),
    @(44427, 3, 6, 374.0648379052369),
    @(44428, 1, 7, 436.4089775561097),
    @(44429, 0, 7, 436.4089775561097),
    @(44430, 0, 6, 374.0648379052369),
    @(44431, 0, 6, 374.0648379052369)
)

$r = $lastRow + 1
foreach ($row in $newData) {
    $destRow = $ws.Range("A$($r):D$($r)")
    # Copy formatting (e.g. the date style/number format used in column A) from the template row
    $template.Copy($destRow)

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]

    $r = $r + 1
}
